# Update "Generate Report for Handback" timestamps.
# These cells hold plain text timestamps (format "yyyy-mm-dd HH:mm:ss"),
# not real Excel date serials, so we assign them as literal strings.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date" for the first file
$wsOverview.Range("G2").Value = "2016-08-28 17:05:17"

# zh-cn!H2 - "Correspond Handoff Datetime" for the first file
$wsZhCn.Range("H2").Value = "2016-08-28 17:05:13"

# zh-cn!K2 - "Correspond Handback DateTime" for the first file
$wsZhCn.Range("K2").Value = "2016-08-28 17:05:29"

# de-de!H2 - "Correspond Handoff Datetime" for the first file (same value as Overview!G2)
$wsDeDe.Range("H2").Value = "2016-08-28 17:05:17"

# de-de!K2 - "Correspond Handback DateTime" for the first file
$wsDeDe.Range("K2").Value = "2016-08-28 17:05:36"
